$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B, shifting the existing B:E data to E:H
$ws.Columns("B:D").Insert()

# New header row values for the newly inserted date columns
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the newly created (empty) B:D cells for each existing analyst row with "UN"
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("B$r").Value = "UN"
    $ws.Range("C$r").Value = "UN"
    $ws.Range("D$r").Value = "UN"
}

# Add two new analyst rows at the bottom of the table
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# Keep the date columns (C:H) at the narrow custom width used throughout the report
$ws.Columns("C:H").ColumnWidth = 7.14

Write-Output "edit complete"
